# Actualización desde MV -datos-
# Adds a new "Agosto.2021" column (BH) to the worksheet, carrying forward
# the last known value (May.2021 / column BG) for each series except the
# first data row, which gets a genuinely updated figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell BH1, copy formatting from BG1 (bold/centered/bordered header)
$ws.Range("BH1").Value = "Agosto.2021"
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for column BH (rows 2-19)
$values = @{
    2  = 53642
    3  = 76549
    4  = 85179
    5  = 92633
    6  = 98915
    7  = 107433
    8  = 100468
    9  = 114106
    10 = 124845
    11 = 133845
    12 = 138668
    13 = 138015
    14 = 141521
    15 = 144093
    16 = 148252
    17 = 154884
    18 = 156422
    19 = 142227
}

foreach ($row in $values.Keys) {
    $ws.Range("BH$row").Value = $values[$row]
}
